# Split the old "Prerequisites" column into Prerequisites / Corequisites /
# Concurrent / Recommended, and pull any "Recommended: ..." suffix out of the
# Prerequisites text into the new "Recommended" column. The old "Terms
# Typically Offered" column (D) moves to G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Insert three new blank columns before the existing "Terms Typically
# Offered" column (D), pushing it to G.
$ws.Columns("D:F").Insert()

# New header row.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

for ($r = 2; $r -le $lastRow; $r++) {
    $prereqCell = $ws.Cells.Item($r, 3)
    $termsCell = $ws.Cells.Item($r, 7)

    $prereq = $prereqCell.Value2
    $terms = $termsCell.Value2

    # Pull a trailing "Recommended: ..." clause out of the Prerequisites
    # text, if present.
    if ($prereq -match "^(.*?)\s*Recommended:\s*(.+)$") {
        $prereqCell.Value = $matches[1]
        $ws.Cells.Item($r, 6).Value = $matches[2]
        $termsCell.Value = "$terms "
    } else {
        $ws.Cells.Item($r, 6).Value = "NA"
    }

    $ws.Cells.Item($r, 4).Value = "NA"
    $ws.Cells.Item($r, 5).Value = "NA"
}
